$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort rows 3-5 ascending by date (Fecha column D), keeping the rest of each
# row's fields together with its own date.
$ws.Range("D3").Value = 44280
$ws.Range("D4").Value = 44284
$ws.Range("D5").Value = 44291

$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 500

$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("P5").Value = 550
